$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New sheet-scoped defined name on Availability (mirrors the existing
#    _FilterDatabase/_FilterDatabase_0/_FilterDatabase_0_0 chain getting one
#    more "_0" suffix appended).
# ---------------------------------------------------------------------------
$wsAvail = $wb.Worksheets.Item("Availability")
$wsAvail.Names.Add("_xlnm._FilterDatabase_0_0_0", "=Availability!`$A`$1:`$S`$16") | Out-Null

# ---------------------------------------------------------------------------
# 2) Availability sheet content fixes: "KO PRODUCT" -> "KO PRODUCTS" in the
#    If-Then value2/value3 columns (H/J), and the Juice row's value3 becomes
#    the tighter "APPLETISER,JUST JUICE" (no space after comma).
# ---------------------------------------------------------------------------
$wsAvail.Range("H11").Value = "KO PRODUCTS"
$wsAvail.Range("J11").Value = "KO PRODUCTS"

$wsAvail.Range("H12").Value = "KO PRODUCTS"

$wsAvail.Range("H13").Value = "KO PRODUCTS"
$wsAvail.Range("J13").Value = "KO PRODUCTS"

$wsAvail.Range("H14").Value = "KO PRODUCTS"
$wsAvail.Range("J14").Value = "KO PRODUCTS"

$wsAvail.Range("H15").Value = "KO PRODUCTS"
$wsAvail.Range("J15").Value = "KO PRODUCTS"

$wsAvail.Range("H16").Value = "KO PRODUCTS"
$wsAvail.Range("J16").Value = "APPLETISER,JUST JUICE"

# ---------------------------------------------------------------------------
# 3) Per-sheet selection / scroll-position bookkeeping (matches the cursor
#    positions left behind after the edit in each sheet).
# ---------------------------------------------------------------------------
$wsKPI = $wb.Worksheets.Item("KPI")
$wsKPI.Activate() | Out-Null
$wsKPI.Range("C22").Select() | Out-Null

$wsPlanogram = $wb.Worksheets.Item("Planogram")
$wsPlanogram.Activate() | Out-Null
$wsPlanogram.Range("C18").Select() | Out-Null

$wsPrice = $wb.Worksheets.Item("Price")
$wsPrice.Activate() | Out-Null
$wsPrice.Range("L8").Select() | Out-Null

$wsSurvey = $wb.Worksheets.Item("Survey")
$wsSurvey.Activate() | Out-Null
$wsSurvey.Range("C17").Select() | Out-Null

$wsSOS = $wb.Worksheets.Item("SOS")
$wsSOS.Activate() | Out-Null
$wsSOS.Range("J9").Select() | Out-Null

$wsCount = $wb.Worksheets.Item("Count")
$wsCount.Activate() | Out-Null
$wsCount.Range("D2").Select() | Out-Null

$wsAvail.Activate() | Out-Null
$wsAvail.Range("G8").Select() | Out-Null
